$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row at position 2 for "Registration Start" and shift the
#    existing schedule rows down.
$ws.Rows(2).Insert()

# The inserted row inherits the header row's (bold/shaded) formatting;
# reset it back to the plain formatting used by the other data rows.
$ws.Range("A2:E2").ClearContents()
$ws.Range("A2:E2").ClearFormats()

$ws.Range("A2").Value = "Sat"
$ws.Range("B2").Value = "10:00am - 10:20am"
$ws.Range("C2").Value = "Hall A"
$ws.Range("D2").Value = "Registration Start"

# 2. Rename "Hall B (Wokshops)" to "Hall B" everywhere on the sheet.
$null = $ws.Cells.Replace("Hall B (Wokshops)", "Hall B")

# 3. Append a new "Closing" row at the end of the schedule.
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Range("A" + $newRow).Value = "Sat"
$ws.Range("B" + $newRow).Value = "06:10pm - 06:30pm"
$ws.Range("C" + $newRow).Value = "Hall A"
$ws.Range("D" + $newRow).Value = "Closing"

Write-Output "done"
